# "commit heure jonas definitif"
# Fills in the Sprint 5 (rows 26-28) journal entries for Jonas Deillon:
#  - reflection text for the sprint (B26)
#  - a new dated task row with hours worked (row 27)
#  - a second task row with hours worked (row 28)
# and tidies up the now-unused blank row 23 that used to belong to the
# previous sprint block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal_DEILLON_JONAS")

# --- New task rows (27-28) --------------------------------------------
$ws.Range("A27").Value = 46010
$ws.Range("B27").Value = "Mise en place su système de récupération de score"
$ws.Range("D27").Value = 5

# --- Sprint reflection text (row 26) ---------------------------------
$ws.Range("B26").Value = "Ce sprint était assez intense, car je me suis beaucoup concentré sur une tâche jusqu’à y parvenir. Ma persévérance a mené à un résultat concluant, avec un système fonctionnel, ce qui constitue un point positif. Ce sprint était principalement axé sur la réalisation, mais nous avons également remarqué que la documentation est primordiale et que nous l’avons quelque peu négligée."

$ws.Range("B28").Value = "Implémentation du l'émulateur en production"
$ws.Range("D28").Value = 1.5

# Row 27 picks up the "continuation row" look (borders/alignment) instead
# of the heavier "first row of block" style it still carried - match it to
# the style already used a couple of rows above (row 21 is a normal,
# non-first row in the previous block) while leaving the date cell (A27)
# untouched.
$ws.Range("B21:D21").Copy() | Out-Null
$ws.Range("B27:D27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 23 is no longer used - remove its contents entirely ----------
$ws.Range("B23:C23").UnMerge() | Out-Null
$ws.Range("B23:D23").Clear() | Out-Null

# --- Update the saved selection / scroll position ----------------------
$ws.Activate()
$ws.Range("I19").Select()
